# Started implementing UV_Subtract function
# Replace the second data row (F1CON/F2CON/F1vF2) with the new sample names,
# and remove the remaining comparison rows that are no longer needed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "27CON"
$ws.Range("B2").Value = "04CON"
$ws.Range("C2").Value = "27v04"

# Remove the now-unused rows 3-7 (F1CON/F3CON/... comparisons) without
# shifting remaining cells, so the sheet dimension shrinks to A1:C2.
$ws.Range("A3:C7").ClearContents()

# Match the author's final selection in the sheet view.
$ws.Range("E8").Select()
